# Added Profit column to spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking text cells as Text so Excel keeps them as
# literal strings instead of coercing them to numbers, then drop the
# temporary format back to Normal so no visible style change remains.
$ws.Range("C2:D6").NumberFormat = "@"

# Header row
$ws.Range("D1").Value = "Profit"

# Data rows
$ws.Range("A2").Value = "Starkonja's Head Silken Hood"
$ws.Range("B2").Value = "2100.0 chaos"
$ws.Range("C2").Value = "3136.875"
$ws.Range("D2").Value = "1036.875"

$ws.Range("A3").Value = "Atziri's Step Slink Boots"
$ws.Range("B3").Value = "50 chaos"
$ws.Range("C3").Value = "755.35"
$ws.Range("D3").Value = "705.35"

$ws.Range("A4").Value = "Devoto's Devotion Nightmare bascinet"
$ws.Range("B4").Value = "420.0 chaos"
$ws.Range("C4").Value = "4606.875"
$ws.Range("D4").Value = "4186.875"

$ws.Range("A5").Value = "Goldrim Leather Cap"
$ws.Range("B5").Value = "80 chaos"
$ws.Range("C5").Value = "145.0"
$ws.Range("D5").Value = "65.0"

$ws.Range("A6").Value = "Greed's Embrace Golden Plate"
$ws.Range("B6").Value = "2100.0 chaos"
$ws.Range("C6").Value = "2765.0"
$ws.Range("D6").Value = "665.0"

# Restore default styling on the cells we temporarily reformatted.
$ws.Range("C2:D6").Style = "Normal"
